$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9.0
$ws.Range("B10").Value = "Sunday, Jan 15"
$ws.Range("C10").Value = "4:10 PM"
$ws.Range("D10").Value = "W92181"
$ws.Range("E10").Value = "London"
$ws.Range("F10").Value = "(LTN)"
$ws.Range("G10").Value = "Wizz Air "
$ws.Range("H10").Value = "A320"
$ws.Range("I10").Value = "(G-WUKF)"
$ws.Range("J10").Value = "4:10 PM"
$ws.Range("K10").Borders.LineStyle = -4142
$ws.Range("L10").Value = "0 hours, 0 minutes"
$ws.Range("M10").Borders.LineStyle = -4142
